$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8138074349, "Toshpulot Toshpulod", 992556550088, "Toyota Camry 1", "14/082025", "2025-08-09 16:20:46", "11:30", "-"),
    @(8138074349, "Yo yo",               992907510905, "Toyota Camry 8", "13/082025", "2025-08-11 09:29:56", "15:00", "-"),
    @(8138074349, "Yo yo",               992907510905, "Toyota Camry 5", "13/082025", "2025-08-11 10:14:06", "11:30", "-"),
    @(8138074349, "Yo yo",               992907510905, "Toyota Camry 5", "15/082025", "2025-08-11 14:56:58", "11:30", "-"),
    @(8138074349, "Yo yo",               992907510905, "Toyota Camry 3", "13/082025", "2025-08-11 15:00:59", "11:30", "-")
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
